$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 54: mark progress 100% for "Agustina" task
$ws.Range("C54").NumberFormat = "0%"
$ws.Range("C54").Value = 1

# Row 55: responsible changes from Lucas to Agustina, progress set to 100%
$ws.Range("B55").Value = "Agustina"
$ws.Range("C55").NumberFormat = "0%"
$ws.Range("C55").Value = 1

# Row 58: mark as "en proceso"
$ws.Range("C58").Value = "en proceso"

# Row 59: assign responsible "Agustina"
$ws.Range("B59").Value = "Agustina"

# Row 61 (new task row): "No asignar cliente en venta de factura B", responsible Lucas
$ws.Range("A61").Value = "No asignar cliente en venta de factura B"
$ws.Range("B61").Value = "Lucas"

# Update selection to reflect new active cell position
$ws.Range("B62").Select()
